$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Update string cells in the precise order needed so the shared-strings
# --- table ends up with the same ordering as the target workbook. ---

# Row 8 labels (Success already present; Failure is new)
$ws1.Range("D8").Value = "Failure"

# Group 1 (columns A/B/C/D) -> Los Angeles Dodgers lefties
$ws1.Range("F3").Value = "Pham"
$ws1.Range("A1").Value = "Los Angeles Dodgers lefties"
$ws1.Range("A3").Value = "Pederson"
$ws1.Range("A4").Value = "Grandal"
$ws1.Range("A5").Value = "Bellinger"
$ws1.Range("A6").Value = "Puig"

# Group 2 (columns F/G/H/I) -> St. Louis Cardinals righties
$ws1.Range("F1").Value = "St. Louis Cardinals righties"
$ws1.Range("F4").Value = "Bader"
$ws1.Range("F5").Value = "Martinez"
$ws1.Range("F6").Value = "Ozuna"

# Group 3 (columns K/L/M/N) -> New York Yankees hitters
$ws1.Range("K1").Value = "New York Yankees hitters"
$ws1.Range("K3").Value = "Gardner"
$ws1.Range("K4").Value = "Hicks"
$ws1.Range("K5").Value = "Walker"
$ws1.Range("K6").Value = "Andujar"

# Remaining "Failure"/"Success" labels
$ws1.Range("I8").Value = "Success"
$ws1.Range("N8").Value = "Failure"

# --- Numeric data updates ---

# Group 1 (Dodgers lefties)
$ws1.Range("B3").Value = 2400
$ws1.Range("C3").Value = 0
$ws1.Range("B4").Value = 3100
$ws1.Range("C4").Value = 0
$ws1.Range("B5").Value = 3400
$ws1.Range("C5").Value = 9
$ws1.Range("B6").Value = 2700
$ws1.Range("C6").Value = 0

# Group 2 (Cardinals righties)
$ws1.Range("G3").Value = 3700
$ws1.Range("H3").Value = 28.2
$ws1.Range("G4").Value = 2800
$ws1.Range("H4").Value = 18.7
$ws1.Range("G5").Value = 3400
$ws1.Range("H5").Value = 3
$ws1.Range("G6").Value = 2800
$ws1.Range("H6").Formula = "=18.7+6.2+6.2"

# Group 3 (Yankees hitters)
$ws1.Range("L3").Value = 3200
$ws1.Range("M3").Value = 15.2
$ws1.Range("L4").Value = 3000
$ws1.Range("M4").Value = 0
$ws1.Range("L5").Value = 2600
$ws1.Range("M5").Value = 3
$ws1.Range("L6").Value = 2700
$ws1.Range("M6").Value = 9.4

# --- Selections (cosmetic, per diff) ---
$ws1.Activate() | Out-Null
$ws1.Range("J10").Select() | Out-Null

$ws2.Activate() | Out-Null
$ws2.Range("A7").Select() | Out-Null

$ws1.Activate() | Out-Null

$wb.Save()
